$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows (2..10, columns B..G) down by one row (into rows 3..11),
# pushing the oldest quarter's row (row 11) out of the table.
for ($r = 10; $r -ge 2; $r--) {
    for ($c = 2; $c -le 7; $c++) {
        $val = $ws.Cells.Item($r, $c).Value2
        $ws.Cells.Item($r + 1, $c).Value = $val
    }
}

# Insert the newest quarter's results into row 2.
$ws.Cells.Item(2, 2).Value = 0.1724578193461484
$ws.Cells.Item(2, 3).Value = 0.39058239716261
$ws.Cells.Item(2, 4).Value = 0.3033305724894426
$ws.Cells.Item(2, 5).Value = 0.550754548314803
$ws.Cells.Item(2, 6).Value = 0.5414156770869448
$ws.Cells.Item(2, 7).Value = 15
